$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - LeetCode #153 Find Minimum in Rotated Sorted Array
$ws.Range("A8").Value = 152
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("B8").Value = "Find Minimum in Rotated Sorted Array"
$ws.Range("B8").Font.Bold = $true
$ws.Range("C8").Value = "Medium"
$ws.Range("C8").Style = "Neutral"
$ws.Range("D8").Value = "Binary Search"
$ws.Range("E8").Value = "O(log n)"
$ws.Range("F8").Value = "Binary search keeping in mind the rotation point in the array. Determine if the middle is part of the `"left`" or `"right`" side of the array."

# Row 9 - LeetCode #167 Two Sum II - Input Array Is Sorted
$ws.Range("A9").Value = 167
$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("B9").Value = "Two Sum II - Input Array Is Sorted"
$ws.Range("B9").Font.Bold = $true
$ws.Range("C9").Value = "Medium"
$ws.Range("C9").Style = "Neutral"
$ws.Range("D9").Value = "Array, Two iters"
$ws.Range("E9").Value = "O(n)"
$ws.Range("F9").Value = "Same idea as a binary search with no middle. Move the right and left iters in increments of 1 based on if the sum of those iters is less or greater than the target."

# Column B width
$ws.Columns("B").ColumnWidth = 32

# Selection update
[void]$ws.Range("F9").Select()

Write-Host "done"
